$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.625.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.59%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.486.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.13%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.613"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.95%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.479.24"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.25%  "

$ws.Range("E9").Value = "  +0.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.183"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.48%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.642"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.82"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.30%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000299"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.84%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.052.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.74%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.486.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.34%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.653.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.43%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.119"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "539.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +13.58%  "

$ws.Range("E22").Value = "  -2.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "93.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.96%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.45%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.71%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "64.55"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.43%  "

$ws.Range("E34").Value = "  -4.79%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "566.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.16%  "

$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "37.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.394"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.82%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0761"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.05%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.33%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.131"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.241.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.87%  "

$ws.Range("E45").Value = "  -3.29%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0435"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.66%  "

$ws.Range("E48").Value = "  -2.26%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "137.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.99%  "
